# Borrar puntos al final de frases en viñetas para versionas Full CV
# Remove trailing periods (and minor whitespace/typo fixes) at the end of
# the supervised-student bullet sentences in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value  = "Maria Paula Moreno Rodríguez (2019 - 2021)"
$ws.Range("E8").Value  = "Andrés Felipe Orozco Serrato (2020 - 2021)"
$ws.Range("E9").Value  = "Danny Ferley Gaitan Rodríguez (2019 - 2020)"
$ws.Range("E10").Value = "Hasbleidy Gamboa Ordoñez (2019 - 2020)"
$ws.Range("E6").Value  = "Andrés Castellanos-Chacón (2017 - 2018; teaching supervision 2019 - Present)"
$ws.Range("E11").Value = "Paula Andrea Betancourt Velandia  (2018 - 2019)"
$ws.Range("E12").Value = "Ana Sofía Gómez Castelblanco (2018 - 2019)"
$ws.Range("E13").Value = "Lina María García Hoyos  (2016 - 2017)"
$ws.Range("E14").Value = "Angie Liliana Pérez Rodríguez  (2016 - 2018)"
$ws.Range("E15").Value = "Lina María Morales Sánchez (2016 - 2017)"
$ws.Range("E16").Value = "Laura Milena Estupiñan Aldana  (2016 - 2017)"
$ws.Range("E17").Value = "Vanesa Díaz Güiza  (2016 - 2018)"
$ws.Range("E18").Value = "Cindy Paola Moncada Gómez (2016 - 2017)"
$ws.Range("E19").Value = "Haydn Ricardo Roldán Morales (2015 - 2016)"
$ws.Range("E20").Value = "Maria Alejandra Abello Mozo  (2017 - 2018)"
$ws.Range("E21").Value = "Natalia Elízabeth Moreno Buitrago (2017 ‑ 2019)"
$ws.Range("E22").Value = "Juan Felipe Pérez Ariza (2017 ‑ 2019)"
$ws.Range("E2").Value  = "Milena Vásquez-Amézquita. Supervised together with  Alicia Salvador"
$ws.Range("E3").Value  = "Francisco Javier Flores. Supervised together with Lisa Chiara Fellin"
$ws.Range("E4").Value  = "Julia Sanz-Vidania. Supervised together with S Craig Roberts"
$ws.Range("E5").Value  = "Adrián Acosta Guerrero. Supervised together with Milena Vásquez-Amézquita"

$ws.Range("E11").Select()
